# Insert a new row for "climate_change_factor_gnrl_hydropower_availability"
# just above the "elasticity_gnrl_rate_occupancy_to_gdppc" row (currently row 4)
# on the first worksheet ("strategy_id-0"), shifting all rows from 4 downward
# one row lower.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4, pushing existing row 4 (and below) to row 5.
$ws.Rows.Item(4).Insert()

# Column A: subsector label
$ws.Cells.Item(4, 1).Value = "General"

# Column B: variable name
$ws.Cells.Item(4, 2).Value = "climate_change_factor_gnrl_hydropower_availability"

# Columns C-G: empty string values (normalize_group, trajgroup_no_vary_q,
# uniform_scaling_q, variable_trajectory_group, variable_trajectory_group_trajectory_type)
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""

# Column H: max_35
$ws.Cells.Item(4, 8).Value = 1

# Column I: min_35
$ws.Cells.Item(4, 9).Value = 0.5

# Columns J (10) through AS (45): yearly trajectory values, all 1
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
